$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used for column D (Price) cells: force the assigned text to
# stay a literal string (the source values are decimal-looking tokens like
# "304.34" or thousand-dotted tokens like "43.374.21") instead of being
# auto-coerced into a floating point number by Excels input parser, then
# drop back to the default "Normal" style so no stray number format sticks
# to the cell (matches the original un-styled inline-string cells).

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Row-by-row Price (D) / Volume 1h (E) updates ---
Set-TextValue $ws.Range('D2') '43.374.21'
$ws.Range('E2').Value = '  +0.01%  '
Set-TextValue $ws.Range('D3') '2.330.61'
$ws.Range('E3').Value = '  -1.24%  '
Set-TextValue $ws.Range('D5') '304.34'
$ws.Range('E5').Value = '  -1.95%  '
Set-TextValue $ws.Range('D6') '100.77'
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E7').Value = '  -3.29%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -3.10%  '
Set-TextValue $ws.Range('D10') '35.19'
$ws.Range('E10').Value = '  -2.67%  '
Set-TextValue $ws.Range('D11') '0.0798'
$ws.Range('E11').Value = '  -2.04%  '
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('E13').Value = '  -3.24%  '
Set-TextValue $ws.Range('D14') '2.689.60'
$ws.Range('E14').Value = '  -1.35%  '
Set-TextValue $ws.Range('D15') '15.67'
$ws.Range('E15').Value = '  -0.16%  '
Set-TextValue $ws.Range('D16') '2.322.32'
$ws.Range('E16').Value = '  -1.16%  '
Set-TextValue $ws.Range('D17') '0.805'
$ws.Range('E17').Value = '  -1.11%  '
Set-TextValue $ws.Range('D18') '43.298.01'
$ws.Range('E18').Value = '  -0.17%  '
Set-TextValue $ws.Range('D19') '11.85'
$ws.Range('E19').Value = '  -1.34%  '
Set-TextValue $ws.Range('D20') '0.0₃0909'
$ws.Range('E20').Value = '  -2.22%  '
Set-TextValue $ws.Range('D21') '6.09'
$ws.Range('E21').Value = '  -2.77%  '
Set-TextValue $ws.Range('D22') '68.17'
$ws.Range('E22').Value = '  -0.28%  '
Set-TextValue $ws.Range('D23') '237.38'
$ws.Range('E23').Value = '  -2.34%  '
Set-TextValue $ws.Range('D24') '1.99'
$ws.Range('E24').Value = '  -3.41%  '
$ws.Range('E25').Value = '  -3.90%  '
$ws.Range('E26').Value = '  -0.28%  '
Set-TextValue $ws.Range('D27') '24.96'
$ws.Range('E27').Value = '  -3.99%  '
$ws.Range('E28').Value = '  -5.73%  '
Set-TextValue $ws.Range('D29') '34.48'
$ws.Range('E29').Value = '  -5.97%  '
Set-TextValue $ws.Range('D30') '165.64'
$ws.Range('E30').Value = '  +1.95%  '
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('E32').Value = '  +0.05%  '
Set-TextValue $ws.Range('D33') '5.06'
$ws.Range('E34').Value = '  -1.98%  '
$ws.Range('E35').Value = '  -4.66%  '
$ws.Range('E38').Value = '  -6.63%  '
Set-TextValue $ws.Range('D39') '1.82'
$ws.Range('E39').Value = '  -6.35%  '
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('E41').Value = '  -3.30%  '
Set-TextValue $ws.Range('D42') '2.43'
$ws.Range('E42').Value = '  -1.56%  '
Set-TextValue $ws.Range('D43') '1.974.48'
$ws.Range('E43').Value = '  -1.16%  '
Set-TextValue $ws.Range('D44') '0.0283'
$ws.Range('E44').Value = '  -2.97%  '
Set-TextValue $ws.Range('D45') '18.55'
$ws.Range('E45').Value = '  -6.21%  '
Set-TextValue $ws.Range('D46') '9.99'
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('E47').Value = '  -5.78%  '
Set-TextValue $ws.Range('D48') '55.89'
$ws.Range('E48').Value = '  -4.47%  '
Set-TextValue $ws.Range('D49') '4.82'
$ws.Range('E49').Value = '  +2.56%  '
Set-TextValue $ws.Range('D50') '2.553.61'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('E51').Value = '  -2.22%  '

# --- Rows 36/37 swap places: Celestia <-> Hedera (name, link, price, volume) ---
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D36') '0.0704'
$ws.Range('E36').Value = '  -5.01%  '

$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D37') '16.73'
$ws.Range('E37').Value = '  -8.77%  '
